# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook is a "Estado de Cuenta" (account statement) report for
# NIT 9013546861. This update replaces the full list of workers in mora
# (7 workers) with a single, new worker, and updates the related totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Remove the six worker rows that are no longer needed -------------
# Rows 16-22 originally held 7 workers. Only one worker remains after the
# update, so rows 17-22 (the other six) are deleted; row 16 is kept and
# edited in place below. Deleting these rows shifts everything below
# (the signature block) up by six rows automatically.
$ws.Range("A17:J22").EntireRow.Delete() | Out-Null

# --- 2. Update the single remaining worker row (row 16) -------------------
$ws.Range("C16").Value = "1051419213"
$ws.Range("D16").Value = "BECKY ALEXANDRA SARA PEREZ"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 30368

# --- 3. Update the summary fields -----------------------------------------
# VALOR MORA (total overdue amount)
$ws.Range("E11").Value = 30368
# Cant. Trabajadores (worker count)
$ws.Range("C13").Value = 1

# --- 4. Re-fit column D (Nombre Trabajador) now that the longest name ----
#        in the sheet is shorter than before ("BECKY ALEXANDRA SARA PEREZ").
$ws.Columns("D").ColumnWidth = 29.3
